# "fix return data timeout" - update scraped numeric counters that changed
# between two scrape passes of the same Shopee listing.
#
# Sheet "product" (rows 2-13, all the same item):
#   I  (stock)              : "3821" -> "3749"
#   P  (product_info JSON)  : embedded "จำนวนสินค้า" "3821" -> "3749"
#   M  (option stock)       : several per-row numeric updates
# Sheet "shop" (row 2):
#   C2 (follower/score text) "177.5พัน" -> "178พัน"
#   D2 (response rate)       "98%"      -> "97%"
#   G2 (followers)           "99.9พัน"  -> "100.1พัน"
#
# Cells in the source file are stored as plain text (t="str"), even though
# their contents look numeric. Excel normally auto-converts a purely
# numeric-looking string typed into a cell into a real number, so for the
# purely-numeric columns (I, M) and the percentage cell (D2) we pre-format
# the cell as Text ("@") before assigning the value to keep it text, just
# like the source data.

$wb = $excel.ActiveWorkbook
$wsProduct = $wb.Worksheets.Item("product")
$wsShop = $wb.Worksheets.Item("shop")

# --- product sheet: column I (stock) 3821 -> 3749, rows 2-13 ---
$wsProduct.Range("I2:I13").NumberFormat = "@"
for ($r = 2; $r -le 13; $r++) {
    $wsProduct.Range("I$r").Value = "3749"
}

# --- product sheet: column P (product_info JSON), replace embedded count ---
$newInfo = "{`n  `"หมวดหมู่`": `"Shopee>เครื่องใช้ในบ้าน>ห้องครัวและห้องอาหาร>อุปกรณ์เบ็ดเตล็ดในครัวอื่นๆ`",`n  `"จำนวนสินค้า`": `"3749`",`n  `"ส่งจาก`": `"จังหวัดราชบุรี`"`n}"
for ($r = 2; $r -le 13; $r++) {
    $wsProduct.Range("P$r").Value = $newInfo
}

# --- product sheet: column M (option stock), per-row updates ---
$mUpdates = @{ 2 = "522"; 3 = "536"; 4 = "371"; 6 = "287"; 7 = "512"; 8 = "421"; 10 = "636"; 11 = "304"; 13 = "47" }
foreach ($r in $mUpdates.Keys) {
    $wsProduct.Range("M$r").NumberFormat = "@"
    $wsProduct.Range("M$r").Value = $mUpdates[$r]
}

# --- shop sheet: row 2 updates ---
$wsShop.Range("C2").Value = "178พัน"
$wsShop.Range("D2").NumberFormat = "@"
$wsShop.Range("D2").Value = "97%"
$wsShop.Range("G2").Value = "100.1พัน"
